# The document originally contains a typo: "new CCoverageMatrix.h()-yCCoverageMatrix.h"
# (em-dash followed by a stray "y" instead of a "->" arrow). The fix:
#   1. Correct "-y" to "->" (em-dash + greater-than).
#   2. Split the run right after the arrow and re-insert the document's
#      "_GoBack" bookmark there (Word automatically relocates a bookmark
#      when Bookmarks.Add is called again with the same name, so the stale
#      bookmark at the end of the document is removed as a side effect).

$d = $word.ActiveDocument
$emDash = [char]0x2014

# Step 1: fix the typo "<emdash>y" -> "<emdash>>" inside the
# "new CCoverageMatrix.h()" line.
[void]$d.Content.Find.Execute(
    "new CCoverageMatrix.h()" + $emDash + "yCCoverageMatrix.h",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "new CCoverageMatrix.h()" + $emDash + ">CCoverageMatrix.h", 2)

# Step 2: locate the point right after "...()-->" and before the trailing
# "CCoverageMatrix.h" so we can drop the bookmark exactly there.
$rng = $d.Content
[void]$rng.Find.Execute(
    "new CCoverageMatrix.h()" + $emDash + ">",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPoint = $rng.End
$bmRange = $d.Range($splitPoint, $splitPoint)

# Re-adding the "_GoBack" bookmark here moves it from wherever it used to be
# (the end of the document) to this new location, removing the old one.
$d.Bookmarks.Add("_GoBack", $bmRange)
